$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the appointment-date / voice-record / form-show fields on row 2
# from the 30/11/2023 09:00-09:04 run to the new 26/12/2023 09:05-09:09 run.
$ws.Range("N2").Value = "date_range`nAppointment Date : 26/12/2023, Time : [ 09:05 AM to 09:09 AM ]"
$ws.Range("AB2").Value = "26/12/2023"
$ws.Range("AR2").Value = "voice_record_26122023"
$ws.Range("AU2").Value = "formshow_26122023"

# The multi-line text in N2 makes Excel auto-expand row 2's height; restore
# the row to its natural auto-fit height so it stays at the sheet default.
$ws.Rows("2:2").AutoFit()
